# Fill in the Activity Log Sheet for Week 9 (Rick / Richard Dobson)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header block ---
$ws.Range("B2").Value = "Richard Dobson"
$ws.Range("G2").Value = 9

# --- Activity rows (row 4 .. row 9) ---
# Row 4: Review progress (Group)
$ws.Range("A4").Value = "Review progress"
$ws.Range("C4").Value = "G"
$ws.Range("D4").Value = (Get-Date -Year 2019 -Month 9 -Day 23 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E4").Value = (Get-Date -Year 2019 -Month 9 -Day 23 -Hour 9 -Minute 0 -Second 0)
$ws.Range("F4").Value = (Get-Date -Year 2019 -Month 9 -Day 23 -Hour 10 -Minute 0 -Second 0)
$ws.Range("G4").Value = 1

# Row 5: Prepare for next meeting (Group)
$ws.Range("A5").Value = "Prepare for next meeting"
$ws.Range("C5").Value = "G"
$ws.Range("D5").Value = (Get-Date -Year 2019 -Month 9 -Day 24 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E5").Value = (Get-Date -Year 2019 -Month 9 -Day 24 -Hour 9 -Minute 0 -Second 0)
$ws.Range("F5").Value = (Get-Date -Year 2019 -Month 9 -Day 24 -Hour 11 -Minute 0 -Second 0)
$ws.Range("G5").Value = 2

# Row 6: Plan next iteration (Group)
$ws.Range("A6").Value = "Plan next iteration"
$ws.Range("C6").Value = "G"
$ws.Range("D6").Value = (Get-Date -Year 2019 -Month 9 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E6").Value = (Get-Date -Year 2019 -Month 9 -Day 25 -Hour 9 -Minute 0 -Second 0)
$ws.Range("F6").Value = (Get-Date -Year 2019 -Month 9 -Day 25 -Hour 11 -Minute 0 -Second 0)
$ws.Range("G6").Value = 2

# Row 7: Work on next iteration (Individual)
$ws.Range("A7").Value = "Work on next iteration"
$ws.Range("C7").Value = "I"
$ws.Range("D7").Value = (Get-Date -Year 2019 -Month 9 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E7").Value = (Get-Date -Year 2019 -Month 9 -Day 26 -Hour 9 -Minute 0 -Second 0)
$ws.Range("F7").Value = (Get-Date -Year 2019 -Month 9 -Day 26 -Hour 2 -Minute 0 -Second 0)
$ws.Range("H7").Value = 5

# Row 8: Work on next iteration (Individual)
$ws.Range("A8").Value = "Work on next iteration"
$ws.Range("C8").Value = "I"
$ws.Range("D8").Value = (Get-Date -Year 2019 -Month 9 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E8").Value = (Get-Date -Year 2019 -Month 9 -Day 27 -Hour 9 -Minute 0 -Second 0)
$ws.Range("F8").Value = (Get-Date -Year 2019 -Month 9 -Day 27 -Hour 15 -Minute 0 -Second 0)
$ws.Range("H8").Value = 6

# Row 9: Work on next iteration (Individual)
$ws.Range("A9").Value = "Work on next iteration"
$ws.Range("C9").Value = "I"
$ws.Range("D9").Value = (Get-Date -Year 2019 -Month 9 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E9").Value = (Get-Date -Year 2019 -Month 9 -Day 28 -Hour 9 -Minute 0 -Second 0)
$ws.Range("F9").Value = (Get-Date -Year 2019 -Month 9 -Day 28 -Hour 1 -Minute 0 -Second 0)
$ws.Range("H9").Value = 4

# Selection as left by the author before saving
$ws.Range("B2:E2").Select()
